# Move the 6 data rows currently on the "New" sheet (rows 2-7) onto the end
# of the "Previously added" sheet (as rows 134-139, with row 133 becoming a
# duplicated header row, matching row 132's style), then trim "New" back
# down to just its header row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# 1) Duplicate the existing "header-as-data" row (132) into row 133.
$ws1.Range("A132:F132").Copy($ws1.Range("A133"))

# 2) Copy the 6 data rows from "New" (A2:F7) into "Previously added" (A134:F139).
$ws2.Range("A2:F7").Copy($ws1.Range("A134"))

# 3) Re-create the hyperlinks on the newly added rows (column A), pointing at
#    the same target URL as the cell text (matches the existing convention
#    used throughout this workbook).
for ($r = 134; $r -le 139; $r++) {
    $target = $ws1.Cells.Item($r, 1).Value()
    $ws1.Hyperlinks.Add($ws1.Cells.Item($r, 1), $target)
}

# Adding a hyperlink resets the cell style to the default "Hyperlink" style,
# so restore the original (pre-existing) link-cell formatting afterwards.
$ws1.Range("A131").Copy()
$ws1.Range("A134:A139").PasteSpecial(-4122)

# 4) Remove the now-duplicated data from "New", leaving just the header row.
$ws2.Rows("2:7").Delete()
$ws2.Hyperlinks.Delete()
